# Use a different email formulation for the last session (no next session):
# reword the "After this session" / "Before the next session" variants of
# the 5 pairs of texts on the "All" sheet (rows 12, 13, 14, 15, 45, columns G/H)
# to start with "I suggest you" after the session/next-session phrase.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All")

$ws.Cells.Item(12,7).Value = 'People are social creatures. This means that what we do is noticed by others around us, and they can react to it. To help you quit smoking, it can be useful to imagine how people who are important to you will receive your non-smoker future self. This can boost your confidence. The reactions of others, such as friends, colleagues, and family, to your future self can be positive, but they could also be negative. It is good to be prepared for both possibilities. After this session, I suggest you grab a pen and paper and answer these 4 questions: 1) How would the people who are important to you react to the new you, who has quit smoking? 2) How would you feel about their reactions? 3) How would you react if the people who are important to you react positively to the new you? 4) How would you react if the people who are important to you react negatively to the new you?'

$ws.Cells.Item(12,8).Value = 'People are social creatures. This means that what we do is noticed by others around us, and they can react to it. To help you quit smoking, it can be useful to imagine how people who are important to you will receive your non-smoker future self. This can boost your confidence. The reactions of others, such as friends, colleagues, and family, to your future self can be positive, but they could also be negative. It is good to be prepared for both possibilities. Before the next session, I suggest you grab a pen and paper and answer these 4 questions: 1) How would the people who are important to you react to the new you, who has quit smoking? 2) How would you feel about their reactions? 3) How would you react if the people who are important to you react positively to the new you? 4) How would you react if the people who are important to you react negatively to the new you?'

$ws.Cells.Item(13,7).Value = 'As social beings, our actions are observed by those around us, who may react in various ways. To boost your confidence when preparing for becoming more physically active, it can be beneficial to envision how those who are important to you will respond to your future physically active self. This can be your colleagues, friends, family, or neighbors, for example. While their reactions may be positive, they could also be negative. So it is wise to anticipate and prepare for both possible outcomes.  After this session, I suggest you grab a pen and paper and answer these 4 questions: 1) How would the people who are important to you react to the new you, who has become more physically active? 2) How would you feel about their reactions? 3) How would you react if the people who are important to you react positively to the new you? 4) How would you react if the people who are important to you react negatively to the new you?'

$ws.Cells.Item(13,8).Value = 'As social beings, our actions are observed by those around us, who may react in various ways. To boost your confidence when preparing for becoming more physically active, it can be beneficial to envision how those who are important to you will respond to your future physically active self. This can be your colleagues, friends, family, or neighbors, for example. While their reactions may be positive, they could also be negative. So it is wise to anticipate and prepare for both possible outcomes.  Before the next session, I suggest you grab a pen and paper and answer these 4 questions: 1) How would the people who are important to you react to the new you, who has become more physically active? 2) How would you feel about their reactions? 3) How would you react if the people who are important to you react positively to the new you? 4) How would you react if the people who are important to you react negatively to the new you?'

$ws.Cells.Item(14,7).Value = 'To increase your confidence that you will succeed in quitting smoking, it can help to think back to previous successes. Have you ever quit smoking before or reduced the number of times you smoked? Maybe you once only had one instead of two cigarettes after dinner. Every success in quitting smoking counts! After this session, I suggest you take some time to think about you previous successes and make a list on a piece of paper. Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what strengths you have that helped you to achieve these things. Write these strengths down on your list so you do not forget them. You can also hang or place your list somewhere in your home so that you are reminded of your successes and strengths more often. The list shows that you can be proud of yourself.'

$ws.Cells.Item(14,8).Value = 'To increase your confidence that you will succeed in quitting smoking, it can help to think back to previous successes. Have you ever quit smoking before or reduced the number of times you smoked? Maybe you once only had one instead of two cigarettes after dinner. Every success in quitting smoking counts! Before the next session, I suggest you take some time to think about your previous successes and make a list on a piece of paper. Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what strengths you have that helped you to achieve these things. Write these strengths down on your list so you do not forget them. You can also hang or place your list somewhere in your home so that you are reminded of your successes and strengths more often. The list shows that you can be proud of yourself.'

$ws.Cells.Item(15,7).Value = 'Reflecting on past successes can boost your confidence in successfully becoming more physically active. After this session, I suggest you think about when you succeeded in being more physically active. Have you ever taken the bike instead of the car, or taken the stairs instead of the escalator? Every small victory counts! Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what you did that helped you to achieve these things. Write it down on your list so you do not forget it. You can also hang or place your list somewhere in your home so you are reminded of your successes more often. These successes show that you can be proud of yourself.'

$ws.Cells.Item(15,8).Value = 'Reflecting on past successes can boost your confidence in successfully becoming more physically active. Before the next session, I suggest you think about when you succeeded in being more physically active. Have you ever taken the bike instead of the car, or taken the stairs instead of the escalator? Every small victory counts! Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what you did that helped you to achieve these things. Write it down on your list so you do not forget it. You can also hang or place your list somewhere in your home so you are reminded of your successes more often. These successes show that you can be proud of yourself.'

$ws.Cells.Item(45,7).Value = 'To increase your confidence that you will succeed in changing your behavior (e.g., quitting smoking), it can help to think back to previous successes. After this session, I suggest you take a moment to reflect on times when you succeeded in something or felt satisfied. Have you ever won a competition? Did you cook something delicious recently? Or maybe you learned a new language? Grab a pen and paper and make a list of your success moments. No success is too small to write down! Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what you did that helped you to achieve your successes. Write it down on your list so you do not forget it. You can also hang or place your list somewhere in your home so that you are reminded of your successes more often. The list shows that you can be proud of yourself.'

$ws.Cells.Item(45,8).Value = 'To increase your confidence that you will succeed in changing your behavior (e.g., quitting smoking), it can help to think back to previous successes. Before the next session, I suggest you take a moment to reflect on times when you succeeded in something or felt satisfied. Have you ever won a competition? Did you cook something delicious recently? Or maybe you learned a new language? Grab a pen and paper and make a list of your success moments. No success is too small to write down! Take a few minutes to make your list before moving on. Then take a closer look at your list. Try to think about what you did that helped you to achieve your successes. Write it down on your list so you do not forget it. You can also hang or place your list somewhere in your home so that you are reminded of your successes more often. The list shows that you can be proud of yourself.'

# The longer reworded text in rows 13 and 14 causes the wrapped cells to need
# extra vertical space, so bump those two row heights by 15pt each (matching
# the target workbook's saved row heights).
$ws.Rows.Item(13).RowHeight = 240
$ws.Rows.Item(14).RowHeight = 210
